# Scheduled-runner update: refresh market-board derived price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve-profit
# sheets, per the latest pull. Pure data refresh - no formulas in this
# workbook, so cells are written directly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7482

$ws.Range("H72").Value = 7482

$ws.Range("H76").Value = 3003
$ws.Range("I76").Value = 3003
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3003
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2688
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 3003
$ws.Range("I79").Value = 3003
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3003
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -1911
$ws.Range("N79").ClearContents()

$ws.Range("H137").Value = 48877.11
$ws.Range("I137").Value = 120385.71
$ws.Range("K137").Value = 361157.13
$ws.Range("M137").Value = -358607.13

$ws.Range("H138").Value = 2609.28
$ws.Range("I138").Value = 1779.3226
$ws.Range("J138").Value = 2982.1594
$ws.Range("K138").Value = 5337.9678
$ws.Range("L138").Value = 8946.4782
$ws.Range("M138").Value = -197.9678000000004
$ws.Range("N138").Value = -19226.4782

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 94.5

$ws.Range("H43").Value = 11837.333
$ws.Range("J43").Value = 15256.5
$ws.Range("L43").Value = 15256.5
$ws.Range("N43").Value = -15882.5

$ws.Range("H44").Value = 69993
$ws.Range("J44").Value = 69993
$ws.Range("L44").Value = 69993
$ws.Range("N44").Value = -70969

$ws.Range("H45").Value = 2650.1667
$ws.Range("I45").Value = 1674.5
$ws.Range("J45").Value = 3625.8333
$ws.Range("K45").Value = 1674.5
$ws.Range("L45").Value = 3625.8333
$ws.Range("M45").Value = -1297.5
$ws.Range("N45").Value = -4379.8333

$ws.Range("H97").Value = 474.2
$ws.Range("I97").Value = 434.32
$ws.Range("J97").Value = 573.9
$ws.Range("K97").Value = 434.32
$ws.Range("L97").Value = 573.9
$ws.Range("M97").Value = 61.68000000000001
$ws.Range("N97").Value = -1565.9

$ws.Range("H102").Value = 1519.9445
$ws.Range("I102").Value = 1057.3334
$ws.Range("K102").Value = 1057.3334
$ws.Range("M102").Value = 564.6666

$ws.Range("H132").Value = 2824.6592
$ws.Range("I132").Value = 2400
$ws.Range("J132").Value = 3957.0833
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 11871.2499
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -16931.2499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3199.7144
$ws.Range("I105").Value = 2999.6667
$ws.Range("K105").Value = 2999.6667
$ws.Range("M105").Value = -1252.6667

$ws.Range("H107").Value = 1882.1818
$ws.Range("I107").Value = 1394.9231
$ws.Range("J107").Value = 2586
$ws.Range("K107").Value = 1394.9231
$ws.Range("L107").Value = 2586
$ws.Range("M107").Value = 525.0769
$ws.Range("N107").Value = -6426

$ws.Range("H134").Value = 2860792.5
$ws.Range("I134").Value = 4204061
$ws.Range("J134").Value = 6346.375
$ws.Range("K134").Value = 12612183
$ws.Range("L134").Value = 19039.125
$ws.Range("M134").Value = -12609648
$ws.Range("N134").Value = -24109.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4997.75
$ws.Range("I31").Value = 1663
$ws.Range("J31").Value = 5474.143
$ws.Range("K31").Value = 1663
$ws.Range("L31").Value = 5474.143
$ws.Range("M31").Value = -1368
$ws.Range("N31").Value = -6064.143

$ws.Range("H34").Value = 4997.75
$ws.Range("I34").Value = 1663
$ws.Range("J34").Value = 5474.143
$ws.Range("K34").Value = 1663
$ws.Range("L34").Value = 5474.143
$ws.Range("M34").Value = -1461
$ws.Range("N34").Value = -5878.143

$ws.Range("H74").Value = 60157
$ws.Range("J74").Value = 60157
$ws.Range("L74").Value = 60157
$ws.Range("N74").Value = -61905

$ws.Range("H77").Value = 60157
$ws.Range("J77").Value = 60157
$ws.Range("L77").Value = 180471
$ws.Range("N77").Value = -189207

$ws.Range("H105").Value = 1739.5454
$ws.Range("I105").Value = 1406.375
$ws.Range("K105").Value = 1406.375
$ws.Range("M105").Value = 340.625

$ws.Range("H132").Value = 3381.7646
$ws.Range("I132").Value = 3151.08
$ws.Range("K132").Value = 9453.24
$ws.Range("M132").Value = -6923.24

$ws.Range("H134").Value = 1847.4849
$ws.Range("I134").Value = 1730.129
$ws.Range("K134").Value = 5190.387
$ws.Range("M134").Value = -2655.387

$ws.Range("H139").Value = 89970
$ws.Range("J139").Value = 89970
$ws.Range("L139").Value = 89970
$ws.Range("N139").Value = -100250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 374.35294
$ws.Range("J92").Value = 531.2
$ws.Range("L92").Value = 1593.6
$ws.Range("N92").Value = -4089.6

$ws.Range("H131").Value = 1850.2632
$ws.Range("J131").Value = 2136.0715
$ws.Range("L131").Value = 6408.2145
$ws.Range("N131").Value = -16488.2145

$ws.Range("H132").Value = 984.43475
$ws.Range("J132").Value = 952.0476
$ws.Range("L132").Value = 8568.428400000001
$ws.Range("N132").Value = -13628.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3148.2856
$ws.Range("I80").Value = 2376
$ws.Range("J80").Value = 4178
$ws.Range("K80").Value = 2376
$ws.Range("L80").Value = 4178
$ws.Range("M80").Value = -1378
$ws.Range("N80").Value = -6174

$ws.Range("H83").Value = 3148.2856
$ws.Range("I83").Value = 2376
$ws.Range("J83").Value = 4178
$ws.Range("K83").Value = 11880
$ws.Range("L83").Value = 20890
$ws.Range("M83").Value = -6888
$ws.Range("N83").Value = -30874

$ws.Range("H95").Value = 120000
$ws.Range("J95").Value = 120000
$ws.Range("L95").Value = 120000
$ws.Range("N95").Value = -125492

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 1990.8462
$ws.Range("I122").Value = 1641.5
$ws.Range("J122").Value = 2146.111
$ws.Range("K122").Value = 4924.5
$ws.Range("L122").Value = 6438.333
$ws.Range("M122").Value = -2474.5
$ws.Range("N122").Value = -11338.333

$ws.Range("H126").Value = 2914.5
$ws.Range("I126").Value = 2783.6667
$ws.Range("K126").Value = 8351.000100000001
$ws.Range("M126").Value = -5881.000100000001

$ws.Range("H132").Value = 2943.18
$ws.Range("I132").Value = 2784.2563
$ws.Range("K132").Value = 8352.768899999999
$ws.Range("M132").Value = -5822.768899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 33338948
$ws.Range("I40").Value = 41670060
$ws.Range("J40").Value = 14500
$ws.Range("K40").Value = 41670060
$ws.Range("L40").Value = 14500
$ws.Range("M40").Value = -41669924
$ws.Range("N40").Value = -14772

$ws.Range("H46").Value = 4326.4
$ws.Range("I46").Value = 1666
$ws.Range("J46").Value = 4795.8823
$ws.Range("K46").Value = 1666
$ws.Range("L46").Value = 4795.8823
$ws.Range("M46").Value = -1478
$ws.Range("N46").Value = -5171.8823

$ws.Range("H132").Value = 81346.38
$ws.Range("I132").Value = 115078.555
$ws.Range("J132").Value = 5449
$ws.Range("K132").Value = 345235.665
$ws.Range("L132").Value = 16347
$ws.Range("M132").Value = -342705.665
$ws.Range("N132").Value = -21407

$ws.Range("H136").Value = 3427.6667
$ws.Range("I136").Value = 3254.6
$ws.Range("J136").Value = 4293
$ws.Range("K136").Value = 9763.799999999999
$ws.Range("L136").Value = 12879
$ws.Range("M136").Value = -7213.799999999999
$ws.Range("N136").Value = -17979

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 45029.5
$ws.Range("I33").Value = 30019
$ws.Range("J33").Value = 60040
$ws.Range("K33").Value = 30019
$ws.Range("L33").Value = 60040
$ws.Range("M33").Value = -29769
$ws.Range("N33").Value = -60540

$ws.Range("H36").Value = 45029.5
$ws.Range("I36").Value = 30019
$ws.Range("J36").Value = 60040
$ws.Range("K36").Value = 30019
$ws.Range("L36").Value = 60040
$ws.Range("M36").Value = -29769
$ws.Range("N36").Value = -60540

$ws.Range("H62").Value = 3787
$ws.Range("I62").Value = 2937
$ws.Range("K62").Value = 2937
$ws.Range("M62").Value = -2313

$ws.Range("H65").Value = 3787
$ws.Range("I65").Value = 2937
$ws.Range("K65").Value = 14685
$ws.Range("M65").Value = -11565

$ws.Range("H70").Value = 37000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 37000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 37000
$ws.Range("N70").Value = -37630
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 37000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 37000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 37000
$ws.Range("N73").Value = -39184
$ws.Range("M73").ClearContents()

$ws.Range("H81").Value = 3647.5
$ws.Range("I81").Value = 2345.6667
$ws.Range("J81").Value = 5321.2856
$ws.Range("K81").Value = 4691.3334
$ws.Range("L81").Value = 10642.5712
$ws.Range("M81").Value = -3630.3334
$ws.Range("N81").Value = -12764.5712

$ws.Range("H84").Value = 3647.5
$ws.Range("I84").Value = 2345.6667
$ws.Range("J84").Value = 5321.2856
$ws.Range("K84").Value = 23456.667
$ws.Range("L84").Value = 53212.856
$ws.Range("M84").Value = -18152.667
$ws.Range("N84").Value = -63820.856

$ws.Range("H113").Value = 821.3570999999999
$ws.Range("I113").Value = 764.75
$ws.Range("J113").Value = 896.8333
$ws.Range("K113").Value = 2294.25
$ws.Range("L113").Value = 2690.4999
$ws.Range("M113").Value = -124.25
$ws.Range("N113").Value = -7030.4999

$ws.Range("H122").Value = 20837286
$ws.Range("I122").Value = 25644626
$ws.Range("J122").Value = 5477.778
$ws.Range("K122").Value = 76933878
$ws.Range("L122").Value = 16433.334
$ws.Range("M122").Value = -76931428
$ws.Range("N122").Value = -21333.334

$ws.Range("H132").Value = 3658.8
$ws.Range("I132").Value = 3537
$ws.Range("J132").Value = 4146
$ws.Range("K132").Value = 10611
$ws.Range("L132").Value = 12438
$ws.Range("M132").Value = -8081
$ws.Range("N132").Value = -17498
